$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.945.54"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.588.07"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.10"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.97"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.31"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "3.052.77"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "62.863.60"
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").Value = "2.593.57"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.30"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "341.89"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.69"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.20"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").Value = "2.712.71"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.33"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("E29").Value = "  +5.65%  "
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "473.06"
$ws.Range("E32").Value = "  +15.19%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.84"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("E35").Value = "  +4.04%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.404"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.03"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  +3.76%  "
$ws.Range("E41").Value = "  -2.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "158.34"
$ws.Range("E42").Value = "  +4.87%  "
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.35"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.633"
$ws.Range("E45").Value = "  +5.16%  "
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0970"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.33"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("E51").Value = "  +1.09%  "
